# ------------------------------------------------------------------
# Overall: A2 numeric 843 -> text "843"
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")
$ws.Range("A2").Value = "'843"
$ws.Range("A2").Style = "Normal"

# ------------------------------------------------------------------
# County: B2:B68 numeric counts -> text
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("County")
$values = @{2="5"; 3="5"; 4="17"; 5="59"; 6="16"; 7="1"; 8="1"; 9="13"; 10="5"; 11="3"; 12="2"; 13="5"; 14="10"; 15="8"; 16="28"; 17="14"; 18="10"; 19="5"; 20="5"; 21="5"; 22="21"; 23="3"; 24="1"; 25="27"; 26="2"; 27="9"; 28="3"; 29="5"; 30="3"; 31="19"; 32="7"; 33="4"; 34="23"; 35="6"; 36="4"; 37="2"; 38="2"; 39="5"; 40="7"; 41="5"; 42="1"; 43="5"; 44="5"; 45="3"; 46="1"; 47="4"; 48="1"; 49="9"; 50="2"; 51="4"; 52="18"; 53="205"; 54="12"; 55="25"; 56="3"; 57="3"; 58="49"; 59="1"; 60="1"; 61="2"; 62="2"; 63="11"; 64="4"; 65="70"; 66="22"; 67="1"; 68="4"}
foreach ($r in $values.Keys) {
    $ws.Range("B$r").Value = "'" + $values[$r]
    $ws.Range("B$r").Style = "Normal"
}

# Row 69 (Prairie County) becomes an all-zero percentage/dollar row
$ws.Range("B69").Value = "'0.00%"
$ws.Range("B69").Style = "Normal"
$ws.Range("C69").Value = "'$0"
$ws.Range("C69").Style = "Normal"
$ws.Range("D69").Value = "'0.00%"
$ws.Range("D69").Style = "Normal"
$ws.Range("E69").Value = "'0.00%"
$ws.Range("E69").Style = "Normal"
$ws.Range("F69").Value = "'0.00%"
$ws.Range("F69").Style = "Normal"

# New row 70: state-wide Total row
$ws.Range("A70").Value = "Total"
$ws.Range("B70").Value = "'843"
$ws.Range("B70").Style = "Normal"
$ws.Range("C70").Value = "'$1,497,047,399"
$ws.Range("C70").Style = "Normal"
$ws.Range("D70").Value = "'7.13%"
$ws.Range("D70").Style = "Normal"
$ws.Range("E70").Value = "'-20.99%"
$ws.Range("E70").Style = "Normal"
$ws.Range("F70").Value = "'75.21%"
$ws.Range("F70").Style = "Normal"

# ------------------------------------------------------------------
# Congressional District: B2:B6 numeric counts -> text
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Congressional District")
$values = @{2="188"; 3="284"; 4="209"; 5="162"; 6="843"}
foreach ($r in $values.Keys) {
    $ws.Range("B$r").Value = "'" + $values[$r]
    $ws.Range("B$r").Style = "Normal"
}

# ------------------------------------------------------------------
# Size: B2:B8 numeric counts -> text
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Size")
$values = @{2="264"; 3="199"; 4="137"; 5="53"; 6="130"; 7="60"; 8="843"}
foreach ($r in $values.Keys) {
    $ws.Range("B$r").Value = "'" + $values[$r]
    $ws.Range("B$r").Style = "Normal"
}

# ------------------------------------------------------------------
# Subsector: B2:B13 numeric counts -> text
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Subsector")
$values = @{2="51"; 3="64"; 4="27"; 5="79"; 6="26"; 7="313"; 8="6"; 9="65"; 10="20"; 11="179"; 12="13"; 13="843"}
foreach ($r in $values.Keys) {
    $ws.Range("B$r").Value = "'" + $values[$r]
    $ws.Range("B$r").Style = "Normal"
}
